$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header "Save" in H1, matching the styled header row
$ws.Range("H1").Value = "Save"
$ws.Range("H1").Style = $ws.Range("G1").Style

# Add 0 values in H2:H4
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 0
